$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '26.804.68'
$ws.Cells.Item(2, 5).Value = '  -1.41%  '
$ws.Cells.Item(3, 4).Value = '1.859.46'
$ws.Cells.Item(3, 5).Value = '  -0.61%  '
$ws.Cells.Item(4, 4).Value = '1.002'
$ws.Cells.Item(4, 5).Value = '  +0.12%  '
$ws.Cells.Item(5, 4).Value = '305.18'
$ws.Cells.Item(5, 5).Value = '  -0.42%  '
$ws.Cells.Item(6, 4).Value = '1.001'
$ws.Cells.Item(6, 5).Value = '  +0.16%  '
$ws.Cells.Item(7, 4).Value = '0.5072'
$ws.Cells.Item(7, 5).Value = '  -2.41%  '
$ws.Cells.Item(8, 4).Value = '0.3648'
$ws.Cells.Item(8, 5).Value = '  -2.43%  '
$ws.Cells.Item(9, 4).Value = '0.07153'
$ws.Cells.Item(9, 5).Value = '  -0.09%  '
$ws.Cells.Item(10, 4).Value = '0.8871'
$ws.Cells.Item(10, 5).Value = '  -0.69%  '
$ws.Cells.Item(11, 4).Value = '20.48'
$ws.Cells.Item(11, 5).Value = '  -1.43%  '
$ws.Cells.Item(12, 4).Value = '1.884.66'
$ws.Cells.Item(12, 5).Value = '  +0.63%  '
$ws.Cells.Item(13, 4).Value = '0.07479'
$ws.Cells.Item(13, 5).Value = '  -0.79%  '
$ws.Cells.Item(14, 4).Value = '93.84'
$ws.Cells.Item(14, 5).Value = '  +3.64%  '
$ws.Cells.Item(15, 4).Value = '5.206'
$ws.Cells.Item(15, 5).Value = '  -1.93%  '
$ws.Cells.Item(16, 4).Value = '1.002'
$ws.Cells.Item(16, 5).Value = '  +0.09%  '
$ws.Cells.Item(17, 4).Value = '0.000008485'
$ws.Cells.Item(17, 5).Value = '  -0.13%  '
$ws.Cells.Item(18, 4).Value = '14.11'
$ws.Cells.Item(18, 5).Value = '  +0.07%  '
$ws.Cells.Item(19, 4).Value = '0.9994'
$ws.Cells.Item(19, 5).Value = '  -0.09%  '
$ws.Cells.Item(20, 4).Value = '26.857.83'
$ws.Cells.Item(20, 5).Value = '  -1.35%  '
$ws.Cells.Item(21, 4).Value = '4.988'
$ws.Cells.Item(21, 5).Value = '  -0.41%  '
$ws.Cells.Item(22, 4).Value = '2.119.54'
$ws.Cells.Item(22, 5).Value = '  +0.54%  '
$ws.Cells.Item(23, 4).Value = '10.31'
$ws.Cells.Item(23, 5).Value = '  -1.45%  '
$ws.Cells.Item(24, 4).Value = '6.344'
$ws.Cells.Item(24, 5).Value = '  -2.08%  '
$ws.Cells.Item(25, 4).Value = '146.89'
$ws.Cells.Item(25, 5).Value = '  +0.69%  '
$ws.Cells.Item(26, 4).Value = '1.766'
$ws.Cells.Item(26, 5).Value = '  -3.68%  '
$ws.Cells.Item(27, 4).Value = '17.83'
$ws.Cells.Item(27, 5).Value = '  -0.97%  '
$ws.Cells.Item(28, 4).Value = '2.084'
$ws.Cells.Item(28, 5).Value = '  -0.22%  '
$ws.Cells.Item(29, 4).Value = '113.45'
$ws.Cells.Item(29, 5).Value = '  +0.13%  '
$ws.Cells.Item(30, 4).Value = '4.662'
$ws.Cells.Item(30, 5).Value = '  -0.06%  '
$ws.Cells.Item(31, 4).Value = '4.682'
$ws.Cells.Item(31, 5).Value = '  -0.12%  '
$ws.Cells.Item(32, 4).Value = '0.09095'
$ws.Cells.Item(32, 5).Value = '  -1.82%  '
$ws.Cells.Item(33, 4).Value = '0.05023'
$ws.Cells.Item(33, 5).Value = '  -2.30%  '
$ws.Cells.Item(34, 4).Value = '0.7451'
$ws.Cells.Item(34, 5).Value = '  +2.29%  '
$ws.Cells.Item(35, 4).Value = '2.949'
$ws.Cells.Item(35, 5).Value = '  -4.37%  '
$ws.Cells.Item(36, 4).Value = '1.146'
$ws.Cells.Item(36, 5).Value = '  -1.36%  '
$ws.Cells.Item(37, 4).Value = '3.208'
$ws.Cells.Item(37, 5).Value = '  +2.84%  '
$ws.Cells.Item(38, 4).Value = '2.494'
$ws.Cells.Item(38, 5).Value = '  -0.83%  '
$ws.Cells.Item(39, 4).Value = '0.01978'
$ws.Cells.Item(39, 5).Value = '  -2.77%  '
$ws.Cells.Item(40, 4).Value = '0.5522'
$ws.Cells.Item(40, 5).Value = '  +3.80%  '
$ws.Cells.Item(41, 4).Value = '1.071'
$ws.Cells.Item(41, 5).Value = '  -0.36%  '
$ws.Cells.Item(42, 4).Value = '6.549'
$ws.Cells.Item(42, 5).Value = '  +0.18%  '
$ws.Cells.Item(43, 4).Value = '115.56'
$ws.Cells.Item(43, 5).Value = '  -1.00%  '
$ws.Cells.Item(44, 4).Value = '8.533'
$ws.Cells.Item(44, 5).Value = '  +2.12%  '
$ws.Cells.Item(45, 4).Value = '0.1479'
$ws.Cells.Item(45, 5).Value = '  +0.15%  '
$ws.Cells.Item(46, 4).Value = '0.4719'
$ws.Cells.Item(46, 5).Value = '  +1.84%  '
$ws.Cells.Item(47, 4).Value = '1.001'
$ws.Cells.Item(47, 5).Value = '  +0.17%  '
$ws.Cells.Item(48, 4).Value = '10.05'
$ws.Cells.Item(48, 5).Value = '  +0.76%  '
$ws.Cells.Item(49, 4).Value = '37.00'
$ws.Cells.Item(49, 5).Value = '  +0.61%  '
$ws.Cells.Item(50, 4).Value = '1.547'
$ws.Cells.Item(50, 5).Value = '  -1.19%  '
$ws.Cells.Item(51, 4).Value = '62.76'
$ws.Cells.Item(51, 5).Value = '  -1.53%  '
